$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M33").Value = -1.153840000000002
$ws.Range("I33").Value = 230.15384
$ws.Range("H33").Value = 220.28572
$ws.Range("K33").Value = 230.15384
$ws.Range("K62").Value = 3501.3
$ws.Range("H62").Value = 3520.25
$ws.Range("M62").Value = -2877.3
$ws.Range("I62").Value = 3501.3
$ws.Range("I65").Value = 3501.3
$ws.Range("K65").Value = 17506.5
$ws.Range("H65").Value = 3520.25
$ws.Range("M65").Value = -14386.5
$ws.Range("L137").Value = 12064.2357
$ws.Range("N137").Value = -17164.2357
$ws.Range("I137").Value = 1595.2084
$ws.Range("M137").Value = -2235.6252
$ws.Range("K137").Value = 4785.6252
$ws.Range("J137").Value = 4021.4119
$ws.Range("H137").Value = 3017.4656
$ws.Range("J138").Value = 3037074
$ws.Range("N138").Value = -9121502
$ws.Range("L138").Value = 9111222
$ws.Range("M138").Value = -7270.400000000001
$ws.Range("I138").Value = 4136.8
$ws.Range("H138").Value = 2278839.8
$ws.Range("K138").Value = 12410.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("J6").Value = 0
$ws.Range("L9").Value = 40004.5
$ws.Range("J9").Value = 40004.5
$ws.Range("N9").Value = -40344.5
$ws.Range("H9").Value = 40004.5
$ws.Range("H20").Value = 40004.5
$ws.Range("N20").Value = -40544.5
$ws.Range("J20").Value = 40004.5
$ws.Range("L20").Value = 40004.5
$ws.Range("H23").Value = 46254.75
$ws.Range("L23").Value = 42503.5
$ws.Range("N23").Value = -43021.5
$ws.Range("J23").Value = 42503.5
$ws.Range("N32").Value = -30574
$ws.Range("J32").Value = 30000
$ws.Range("M32").Value = -14183.36
$ws.Range("L32").Value = 30000
$ws.Range("I32").Value = 14470.36
$ws.Range("K32").Value = 14470.36
$ws.Range("H32").Value = 17867.469
$ws.Range("N37").Value = -24796
$ws.Range("H37").Value = 18832.5
$ws.Range("L37").Value = 24250
$ws.Range("J37").Value = 24250
$ws.Range("L44").Value = 100043400
$ws.Range("N44").Value = -100044376
$ws.Range("H44").Value = 100043400
$ws.Range("J44").Value = 100043400
$ws.Range("J55").Value = 142904290
$ws.Range("L55").Value = 142904290
$ws.Range("H55").Value = 142904290
$ws.Range("N55").Value = -142904920
$ws.Range("N61").Value = -13418.1875
$ws.Range("H61").Value = 6143.7036
$ws.Range("I61").Value = 3259.2896
$ws.Range("J61").Value = 12994.1875
$ws.Range("M61").Value = -3047.2896
$ws.Range("K61").Value = 3259.2896
$ws.Range("L61").Value = 12994.1875
$ws.Range("H80").Value = 142881630
$ws.Range("N80").Value = -250027976
$ws.Range("J80").Value = 250025980
$ws.Range("L80").Value = 250025980
$ws.Range("J83").Value = 250025980
$ws.Range("N83").Value = -750087924
$ws.Range("H83").Value = 142881630
$ws.Range("L83").Value = 750077940
$ws.Range("K132").Value = 5218.7001
$ws.Range("N132").Value = -24518.201
$ws.Range("J132").Value = 6486.067
$ws.Range("I132").Value = 1739.5667
$ws.Range("M132").Value = -2688.7001
$ws.Range("L132").Value = 19458.201
$ws.Range("H132").Value = 4112.817
$ws.Range("N136").Value = -44082.5625
$ws.Range("J136").Value = 12994.1875
$ws.Range("K136").Value = 9777.8688
$ws.Range("L136").Value = 38982.5625
$ws.Range("M136").Value = -7227.8688
$ws.Range("I136").Value = 3259.2896
$ws.Range("H136").Value = 6143.7036
$ws.Range("L137").Value = 64796.668
$ws.Range("N137").Value = -74996.66800000001
$ws.Range("J137").Value = 64796.668
$ws.Range("H137").Value = 64796.668

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H104").Value = 59999.855
$ws.Range("J104").Value = 59999.855
$ws.Range("N104").Value = -66987.85500000001
$ws.Range("L104").Value = 59999.855
$ws.Range("H116").Value = 79800
$ws.Range("N116").Value = -88978
$ws.Range("L116").Value = 79800
$ws.Range("J116").Value = 79800
$ws.Range("J134").Value = 69977
$ws.Range("I134").Value = 2707.2778
$ws.Range("L134").Value = 209931
$ws.Range("N134").Value = -215001
$ws.Range("K134").Value = 8121.8334
$ws.Range("M134").Value = -5586.8334
$ws.Range("H134").Value = 22492.49

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("K93").Value = 12252.5
$ws.Range("L93").Value = 20000
$ws.Range("I93").Value = 12252.5
$ws.Range("H93").Value = 13802
$ws.Range("M93").Value = -10380.5
$ws.Range("J93").Value = 20000
$ws.Range("N93").Value = -23744
$ws.Range("N94").Value = -1991.75
$ws.Range("I94").Value = 1290.6364
$ws.Range("M94").Value = -839.6364000000001
$ws.Range("L94").Value = 1089.75
$ws.Range("H94").Value = 1185.826
$ws.Range("K94").Value = 1290.6364
$ws.Range("J94").Value = 1089.75
$ws.Range("K132").Value = 5679.7827
$ws.Range("I132").Value = 1893.2609
$ws.Range("M132").Value = -3149.7827
$ws.Range("H132").Value = 2262.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("J2").Value = 53.677418
$ws.Range("H2").Value = 45.102566
$ws.Range("M2").Value = 41.75
$ws.Range("N2").Value = -548.064508
$ws.Range("I2").Value = 11.875
$ws.Range("L2").Value = 322.064508
$ws.Range("K2").Value = 71.25
$ws.Range("H3").Value = 6680.1787
$ws.Range("N3").Value = -30114.386
$ws.Range("K3").Value = 11504.0001
$ws.Range("J3").Value = 9963.462
$ws.Range("M3").Value = -11392.0001
$ws.Range("L3").Value = 29890.386
$ws.Range("I3").Value = 3834.6667
$ws.Range("J34").Value = 4264.516
$ws.Range("I34").Value = 130
$ws.Range("M34").Value = -306
$ws.Range("H34").Value = 3594.054
$ws.Range("N34").Value = -12961.548
$ws.Range("L34").Value = 12793.548
$ws.Range("K34").Value = 390
$ws.Range("J39").Value = 8751.678
$ws.Range("N39").Value = -26843.034
$ws.Range("L39").Value = 26255.034
$ws.Range("H39").Value = 8751.678
$ws.Range("M55").ClearContents()
$ws.Range("J55").Value = 8666.5
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 25999.5
$ws.Range("H55").Value = 8666.5
$ws.Range("N55").Value = -26353.5
$ws.Range("J115").Value = 2507.5715
$ws.Range("M115").Value = -5209.75
$ws.Range("L115").Value = 7522.7145
$ws.Range("H115").Value = 2369.6365
$ws.Range("K115").Value = 6384.75
$ws.Range("I115").Value = 2128.25
$ws.Range("N115").Value = -9872.7145
$ws.Range("I122").Value = 650.3077
$ws.Range("L122").Value = 15634.125
$ws.Range("M122").Value = -3402.7693
$ws.Range("N122").Value = -20534.125
$ws.Range("H122").Value = 1249.931
$ws.Range("J122").Value = 1737.125
$ws.Range("K122").Value = 5852.7693
$ws.Range("H123").Value = 6999.857
$ws.Range("N123").Value = -28400.0005
$ws.Range("K123").Value = 5997
$ws.Range("I123").Value = 1999
$ws.Range("L123").Value = 23500.0005
$ws.Range("M123").Value = -3547
$ws.Range("J123").Value = 7833.3335
$ws.Range("K132").Value = 13081.5
$ws.Range("I132").Value = 1453.5
$ws.Range("M132").Value = -10551.5
$ws.Range("H132").Value = 1797.4445
$ws.Range("H140").Value = 1802.4286
$ws.Range("M140").Value = 1324.6925
$ws.Range("K140").Value = 3855.3075
$ws.Range("I140").Value = 1285.1025

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N126").Value = -15895.769
$ws.Range("I126").Value = 1912.4445
$ws.Range("M126").Value = -3267.333500000001
$ws.Range("K126").Value = 5737.333500000001
$ws.Range("J126").Value = 3651.923
$ws.Range("L126").Value = 10955.769
$ws.Range("H126").Value = 2940.318
$ws.Range("K132").Value = 17044.6671
$ws.Range("I132").Value = 5681.5557
$ws.Range("M132").Value = -14514.6671
$ws.Range("H132").Value = 6193.795
$ws.Range("L135").Value = 49911.43
$ws.Range("N135").Value = -60051.43
$ws.Range("H135").Value = 49911.43
$ws.Range("J135").Value = 49911.43

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K16").Value = 615.5333000000001
$ws.Range("M16").Value = -445.5333000000001
$ws.Range("N16").ClearContents()
$ws.Range("I16").Value = 615.5333000000001
$ws.Range("L16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("H16").Value = 615.5333000000001
$ws.Range("I25").Value = 9999.666999999999
$ws.Range("H25").Value = 10998.167
$ws.Range("J25").Value = 11996.667
$ws.Range("K25").Value = 9999.666999999999
$ws.Range("N25").Value = -12456.667
$ws.Range("L25").Value = 11996.667
$ws.Range("M25").Value = -9769.666999999999
$ws.Range("L46").Value = 490
$ws.Range("H46").Value = 506.66666
$ws.Range("J46").Value = 490
$ws.Range("N46").Value = -866
$ws.Range("M46").Value = -327
$ws.Range("I46").Value = 515
$ws.Range("K46").Value = 515
$ws.Range("I82").Value = 1120.1111
$ws.Range("K82").Value = 1120.1111
$ws.Range("N82").Value = -4944.2
$ws.Range("H82").Value = 1794.4783
$ws.Range("L82").Value = 4222.2
$ws.Range("M82").Value = -759.1111000000001
$ws.Range("J82").Value = 4222.2
$ws.Range("N85").Value = -6718.2
$ws.Range("L85").Value = 4222.2
$ws.Range("M85").Value = 127.8888999999999
$ws.Range("I85").Value = 1120.1111
$ws.Range("K85").Value = 1120.1111
$ws.Range("H85").Value = 1794.4783
$ws.Range("J85").Value = 4222.2
$ws.Range("K132").Value = 19909.125
$ws.Range("N132").Value = -15362.9999
$ws.Range("J132").Value = 3434.3333
$ws.Range("I132").Value = 6636.375
$ws.Range("M132").Value = -17379.125
$ws.Range("L132").Value = 10302.9999
$ws.Range("H132").Value = 5763.091
$ws.Range("J134").Value = 65427
$ws.Range("L134").Value = 65427
$ws.Range("N134").Value = -75567
$ws.Range("H134").Value = 65427

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K132").Value = 3723.1248
$ws.Range("N132").Value = -12268.6469
$ws.Range("J132").Value = 2402.8823
$ws.Range("I132").Value = 1241.0416
$ws.Range("M132").Value = -1193.1248
$ws.Range("L132").Value = 7208.646900000001
$ws.Range("H132").Value = 1722.7805
